$d = $word.ActiveDocument

$replacements = @(
    @("659×7=", "456×6="),
    @("292×4=", "651×9="),
    @("616×9=", "540×7="),
    @("296×8=", "646×8="),
    @("850×4=", "132×4="),
    @("946×7=", "905×4="),
    @("566×5=", "979×3="),
    @("818×7=", "197×8="),
    @("486×2=", "193×2="),
    @("502×8=", "225×8="),
    @("542×2=", "902×2="),
    @("389×4=", "685×8="),
    @("502×3=", "231×8="),
    @("446×6=", "180×4="),
    @("515×5=", "739×2="),
    @("774×8=", "326×2="),
    @("774×3=", "255×3="),
    @("862×8=", "274×8="),
    @("981×8=", "709×7="),
    @("129×9=", "867×6="),
    @("724×7=", "848×4="),
    @("964×7=", "114×7="),
    @("459×9=", "425×5="),
    @("656×7=", "586×4="),
    @("718×2=", "864×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $new, 2)
}
